$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 64: Forged from the Void | Void Glue
$ws.Cells.Item(64, 8).Value = 2985.524
$ws.Cells.Item(64, 9).Value = 2717.818
$ws.Cells.Item(64, 10).Value = 3280
$ws.Cells.Item(64, 11).Value = 2717.818
$ws.Cells.Item(64, 12).Value = 3280
$ws.Cells.Item(64, 13).Value = -2469.818
$ws.Cells.Item(64, 14).Value = -3776

# Row 67: Dodging the Draft (L) | Void Glue
$ws.Cells.Item(67, 8).Value = 2985.524
$ws.Cells.Item(67, 9).Value = 2717.818
$ws.Cells.Item(67, 10).Value = 3280
$ws.Cells.Item(67, 11).Value = 2717.818
$ws.Cells.Item(67, 12).Value = 3280
$ws.Cells.Item(67, 13).Value = -1859.818
$ws.Cells.Item(67, 14).Value = -4996

# Row 128: Nearly There | Kumbhiraskin Grimoire
$ws.Cells.Item(128, 8).Value = 90000
$ws.Cells.Item(128, 10).Value = 90000
$ws.Cells.Item(128, 12).Value = 90000
$ws.Cells.Item(128, 14).Value = -99960

# Row 137: Cutting Edge of Culinary Quality | Magnesia Whetstone
$ws.Cells.Item(137, 8).Value = 2300.0908
$ws.Cells.Item(137, 9).Value = 1643
$ws.Cells.Item(137, 11).Value = 4929
$ws.Cells.Item(137, 13).Value = -2379

# Row 138: All-night Crafting | Cunning Craftsman's Tisane
$ws.Cells.Item(138, 8).Value = 4025.4744
$ws.Cells.Item(138, 9).Value = 1879.4
$ws.Cells.Item(138, 10).Value = 5366.771
$ws.Cells.Item(138, 11).Value = 5638.200000000001
$ws.Cells.Item(138, 12).Value = 16100.313
$ws.Cells.Item(138, 13).Value = -498.2000000000007
$ws.Cells.Item(138, 14).Value = -26380.313

$ws = $wb.Worksheets.Item("ARM")
# Row 24: A Firm Hand | Iron Gauntlets
$ws.Cells.Item(24, 8).Value = 21000
$ws.Cells.Item(24, 10).Value = 21000
$ws.Cells.Item(24, 12).Value = 21000
$ws.Cells.Item(24, 14).Value = -21748

# Row 32: Ingot We Trust | Steel Ingot
$ws.Cells.Item(32, 8).Value = 32288674
$ws.Cells.Item(32, 9).Value = 62515396
$ws.Cells.Item(32, 10).Value = 46839.6
$ws.Cells.Item(32, 11).Value = 62515396
$ws.Cells.Item(32, 12).Value = 46839.6
$ws.Cells.Item(32, 13).Value = -62515109
$ws.Cells.Item(32, 14).Value = -47413.6

# Row 100: En Garde and on Guard | Doman Iron Gauntlets of Fending
$ws.Cells.Item(100, 8).Value = 21000
$ws.Cells.Item(100, 10).Value = 21000
$ws.Cells.Item(100, 12).Value = 21000
$ws.Cells.Item(100, 14).Value = -23164

# Row 128: Heading toward Bankruptcy | Manganese Helm of the Falling Dragon
$ws.Cells.Item(128, 8).Value = 38429
$ws.Cells.Item(128, 10).Value = 38429
$ws.Cells.Item(128, 12).Value = 38429
$ws.Cells.Item(128, 14).Value = -48389

# Row 132: Don't Bore Me, Ore Me | Mountain Chromite Ingot
$ws.Cells.Item(132, 8).Value = 20003030
$ws.Cells.Item(132, 9).Value = 26316558
$ws.Cells.Item(132, 10).Value = 10188.167
$ws.Cells.Item(132, 11).Value = 78949674
$ws.Cells.Item(132, 12).Value = 30564.501
$ws.Cells.Item(132, 13).Value = -78947144
$ws.Cells.Item(132, 14).Value = -35624.501

# Row 134: Brace for More Vambraces | Ruthenium Vambraces of Maiming
$ws.Cells.Item(134, 8).Value = 47000
$ws.Cells.Item(134, 10).Value = 47000
$ws.Cells.Item(134, 12).Value = 47000
$ws.Cells.Item(134, 14).Value = -57140

$ws = $wb.Worksheets.Item("BSM")
# Row 134: Ruthenium Supremium | Ruthenium Ingot
$ws.Cells.Item(134, 8).Value = 2102196.2
$ws.Cells.Item(134, 9).Value = 4973.0938
$ws.Cells.Item(134, 10).Value = 5297964.5
$ws.Cells.Item(134, 11).Value = 14919.2814
$ws.Cells.Item(134, 12).Value = 15893893.5
$ws.Cells.Item(134, 13).Value = -12384.2814
$ws.Cells.Item(134, 14).Value = -15898963.5

# Row 135: Axes to the Maxes | Ruthenium War Axe
$ws.Cells.Item(135, 8).Value = 98000
$ws.Cells.Item(135, 10).Value = 98000
$ws.Cells.Item(135, 12).Value = 98000
$ws.Cells.Item(135, 14).Value = -108140

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found | Walnut Lumber
$ws.Cells.Item(31, 8).Value = 2570.4546
$ws.Cells.Item(31, 9).Value = 2166.3044
$ws.Cells.Item(31, 10).Value = 3500
$ws.Cells.Item(31, 11).Value = 2166.3044
$ws.Cells.Item(31, 12).Value = 3500
$ws.Cells.Item(31, 13).Value = -1871.3044
$ws.Cells.Item(31, 14).Value = -4090

# Row 34: Armoires of the Rich and Famous | Walnut Lumber
$ws.Cells.Item(34, 8).Value = 2570.4546
$ws.Cells.Item(34, 9).Value = 2166.3044
$ws.Cells.Item(34, 10).Value = 3500
$ws.Cells.Item(34, 11).Value = 2166.3044
$ws.Cells.Item(34, 12).Value = 3500
$ws.Cells.Item(34, 13).Value = -1964.3044
$ws.Cells.Item(34, 14).Value = -3904

# Row 134: Wood You Be Quiet | Ceiba Lumber
$ws.Cells.Item(134, 8).Value = 1264.7213
$ws.Cells.Item(134, 9).Value = 1218.7291
$ws.Cells.Item(134, 10).Value = 1434.5385
$ws.Cells.Item(134, 11).Value = 3656.1873
$ws.Cells.Item(134, 12).Value = 4303.6155
$ws.Cells.Item(134, 13).Value = -1121.1873
$ws.Cells.Item(134, 14).Value = -9373.6155

# Row 141: No Greater Treasure | Claro Walnut Necklace of Gathering
$ws.Cells.Item(141, 8).Value = 58000
$ws.Cells.Item(141, 9).Value = 16000
$ws.Cells.Item(141, 10).Value = 65000
$ws.Cells.Item(141, 11).Value = 16000
$ws.Cells.Item(141, 12).Value = 65000
$ws.Cells.Item(141, 13).Value = -10820
$ws.Cells.Item(141, 14).Value = -75360

$ws = $wb.Worksheets.Item("CUL")
# Row 47: Winter of Our Discontent | Mugwort Carp
$ws.Cells.Item(47, 8).Value = 456.25
$ws.Cells.Item(47, 9).Value = 12.5
$ws.Cells.Item(47, 10).Value = 900
$ws.Cells.Item(47, 11).Value = 37.5
$ws.Cells.Item(47, 12).Value = 2700
$ws.Cells.Item(47, 13).Value = 393.5
$ws.Cells.Item(47, 14).Value = -3562

# Row 60: Drinking to Your Health | Mulled Tea
$ws.Cells.Item(60, 8).Value = 1812.5
$ws.Cells.Item(60, 9).Value = 500
$ws.Cells.Item(60, 11).Value = 1500
$ws.Cells.Item(60, 13).Value = -1249

# Row 109: Cure for What Ails | Purple Carrot Juice
$ws.Cells.Item(109, 8).Value = 3317.7058
$ws.Cells.Item(109, 9).Value = 794.2857
$ws.Cells.Item(109, 10).Value = 5084.1
$ws.Cells.Item(109, 11).Value = 2382.8571
$ws.Cells.Item(109, 12).Value = 15252.3
$ws.Cells.Item(109, 13).Value = -1342.8571
$ws.Cells.Item(109, 14).Value = -17332.3

# Row 131: The Mountain Steeped | Tsai tou Vounou
$ws.Cells.Item(131, 8).Value = 872.14
$ws.Cells.Item(131, 9).Value = 556.6667
$ws.Cells.Item(131, 10).Value = 881.8969
$ws.Cells.Item(131, 11).Value = 1670.0001
$ws.Cells.Item(131, 12).Value = 2645.6907
$ws.Cells.Item(131, 13).Value = 3369.9999
$ws.Cells.Item(131, 14).Value = -12725.6907

$ws = $wb.Worksheets.Item("GSM")
# Row 97: If I'd a Koppranickel for Every Time... | Koppranickel Ingot
$ws.Cells.Item(97, 8).Value = 835.9231
$ws.Cells.Item(97, 9).Value = 851.5454999999999
$ws.Cells.Item(97, 10).Value = 750
$ws.Cells.Item(97, 11).Value = 851.5454999999999
$ws.Cells.Item(97, 12).Value = 750
$ws.Cells.Item(97, 13).Value = -355.5454999999999
$ws.Cells.Item(97, 14).Value = -1742

$ws = $wb.Worksheets.Item("LTW")
# Row 16: Saddle Sore | Hard Leather
$ws.Cells.Item(16, 8).Value = 59525310
$ws.Cells.Item(16, 9).Value = 5102960.5
$ws.Cells.Item(16, 10).Value = 250003550
$ws.Cells.Item(16, 11).Value = 5102960.5
$ws.Cells.Item(16, 12).Value = 250003550
$ws.Cells.Item(16, 13).Value = -5102790.5
$ws.Cells.Item(16, 14).Value = -250003890

# Row 22: Skin off Their Backs | Aldgoat Leather
$ws.Cells.Item(22, 8).Value = 12658227
$ws.Cells.Item(22, 9).Value = 12658227
$ws.Cells.Item(22, 11).Value = 12658227
$ws.Cells.Item(22, 13).Value = -12657932

# Row 27: Fire and Hide | Aldgoat Leather
$ws.Cells.Item(27, 8).Value = 12658227
$ws.Cells.Item(27, 9).Value = 12658227
$ws.Cells.Item(27, 11).Value = 12658227
$ws.Cells.Item(27, 13).Value = -12658120

# Row 46: Supply Side Logic | Boar Leather
$ws.Cells.Item(46, 8).Value = 5952837.5
$ws.Cells.Item(46, 9).Value = 6944894
$ws.Cells.Item(46, 10).Value = 500
$ws.Cells.Item(46, 11).Value = 6944894
$ws.Cells.Item(46, 12).Value = 500
$ws.Cells.Item(46, 13).Value = -6944706
$ws.Cells.Item(46, 14).Value = -876

# Row 132: Tenets of Tanning | Silver Lobo Leather
$ws.Cells.Item(132, 8).Value = 21745240
$ws.Cells.Item(132, 9).Value = 52634170
$ws.Cells.Item(132, 10).Value = 8585.296
$ws.Cells.Item(132, 11).Value = 157902510
$ws.Cells.Item(132, 12).Value = 25755.888
$ws.Cells.Item(132, 13).Value = -157899980
$ws.Cells.Item(132, 14).Value = -30815.888

# Row 136: Respect for Br'aax | Br'aax Leather
$ws.Cells.Item(136, 8).Value = 2964.195
$ws.Cells.Item(136, 9).Value = 2384.4666
$ws.Cells.Item(136, 10).Value = 4545.273
$ws.Cells.Item(136, 11).Value = 7153.399800000001
$ws.Cells.Item(136, 12).Value = 13635.819
$ws.Cells.Item(136, 13).Value = -4603.399800000001
$ws.Cells.Item(136, 14).Value = -18735.819

$ws = $wb.Worksheets.Item("WVR")
# Row 44: Edmelle's Hair | Linen Wedge Cap of Gathering
$ws.Cells.Item(44, 8).Value = 5500
$ws.Cells.Item(44, 9).Value = 5000
$ws.Cells.Item(44, 11).Value = 5000
$ws.Cells.Item(44, 13).Value = -4446

# Row 92: Modest Beginnings | Bloodhempen Culottes of Casting
$ws.Cells.Item(92, 8).Value = 19500
$ws.Cells.Item(92, 10).Value = 19500
$ws.Cells.Item(92, 12).Value = 19500
$ws.Cells.Item(92, 14).Value = -24492

# Row 126: A Polished Purchase | Snow Linen
$ws.Cells.Item(126, 8).Value = 1662
$ws.Cells.Item(126, 9).Value = 1060.6
$ws.Cells.Item(126, 10).Value = 3666.6667
$ws.Cells.Item(126, 11).Value = 3181.8
$ws.Cells.Item(126, 12).Value = 11000.0001
$ws.Cells.Item(126, 13).Value = -711.7999999999997
$ws.Cells.Item(126, 14).Value = -15940.0001

# Row 136: Weaving the Envelope | Sarcenet Cloth
$ws.Cells.Item(136, 8).Value = 3414.4314
$ws.Cells.Item(136, 9).Value = 5303.625
$ws.Cells.Item(136, 10).Value = 1735.1482
$ws.Cells.Item(136, 11).Value = 15910.875
$ws.Cells.Item(136, 12).Value = 5205.444600000001
$ws.Cells.Item(136, 13).Value = -13360.875
$ws.Cells.Item(136, 14).Value = -10305.4446
